$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 2 and row 3 and need to be swapped.
$cols = @("A", "B", "E", "F", "G", "H", "I", "Y", "AA", "AI")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $v2 = $ws.Range($addr2).Value()
    $v3 = $ws.Range($addr3).Value()
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# L2 becomes an empty cell (present but blank) while L3 (previously an empty
# cell) is cleared entirely, matching the row swap.
$ws.Range("L2").Value = ""
$ws.Range("L3").ClearContents()
